$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The merged "Budget Précédent" / "Budget\nActuel" column headers were split
# into the shorter "Précédent" / "Actuel" labels.
$ws.Range("B1").Value = "Précédent"
$ws.Range("C1").Value = "Actuel"

# Fix the Voyage note typo: "1,6 000 $" -> "16 000 $"
$ws.Range("G14").Value = "10 000 $ pour effectuer des inspections, 16 000 $ pour les déplacements généraux"

# Leave the selection on the cell that was last edited, as in the saved file.
$ws.Range("G14").Select()
